# Update "Glosario Ranking Muni_2018_GT.xlsx"
#
# The "Código" column (F) on Hoja1 previously stored bare codes like
# "04-01", "01-01", "00-00", etc. This adds a "COD " prefix to every code
# in the data rows (F2:F31) while leaving the column header in F1
# ("Código ") untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- F2:F31 : prefix each existing code with "COD " ------------------
# A leading apostrophe forces a text/literal assignment so the cell keeps
# its original number format / quote-prefix styling (s="6" / s="7")
# instead of Excel re-stamping it with a fresh "General" style.
for ($r = 2; $r -le 31; $r++) {
    $cell = $ws.Cells.Item($r, 6)   # column F
    $oldCode = $cell.Value2
    $cell.Value2 = "'COD " + $oldCode
}

# --- view state: scroll the sheet over one column and move the
#     selection from F2 to H4, matching the author's saved view -------
$win = $excel.ActiveWindow
$win.ScrollColumn = 5
$win.ScrollRow = 1
$ws.Range("H4").Select() | Out-Null

# --- restore the workbook window to a maximized-looking geometry -----
$wbWin = $wb.Windows.Item(1)
$wbWin.Left = -120
$wbWin.Top = -120
$wbWin.Width = 20730
$wbWin.Height = 11160

Write-Host "Updated F2:F31 codes with 'COD ' prefix; selection -> H4"
